# Update state_2005_2006.xlsx to the September/October version:
# refresh the tau / lambda / lambda_hat (and T / T_hat) regression
# output values on each of the three sheets (HSV Log OLS, HSV PPML,
# HSVT NLLSQ) for rows 2-52.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("B2").Value = -0.014
$ws1.Range("C2").Value = 0.777
$ws1.Range("D2").Value = 0.916
$ws1.Range("B3").Value = 0.047
$ws1.Range("C3").Value = 1.63
$ws1.Range("D3").Value = 0.939
$ws1.Range("B4").Value = -0.003
$ws1.Range("C4").Value = 0.9
$ws1.Range("D4").Value = 0.932
$ws1.Range("B5").Value = -0.004
$ws1.Range("C5").Value = 0.88
$ws1.Range("D5").Value = 0.922
$ws1.Range("B6").Value = 0.015
$ws1.Range("C6").Value = 1.114
$ws1.Range("B7").Value = 0.0
$ws1.Range("C7").Value = 0.925
$ws1.Range("D7").Value = 0.925
$ws1.Range("B8").Value = 0.011
$ws1.Range("C8").Value = 1.037
$ws1.Range("D8").Value = 0.911
$ws1.Range("B9").Value = 0.018
$ws1.Range("C9").Value = 1.137
$ws1.Range("D9").Value = 0.922
$ws1.Range("B10").Value = 0.015
$ws1.Range("C10").Value = 1.085
$ws1.Range("D10").Value = 0.91
$ws1.Range("B11").Value = -0.016
$ws1.Range("C11").Value = 0.775
$ws1.Range("D11").Value = 0.937
$ws1.Range("B12").Value = 0.006
$ws1.Range("C12").Value = 0.992
$ws1.Range("B13").Value = 0.012
$ws1.Range("C13").Value = 1.058
$ws1.Range("D13").Value = 0.919
$ws1.Range("B14").Value = 0.012
$ws1.Range("C14").Value = 1.07
$ws1.Range("D14").Value = 0.929
$ws1.Range("B15").Value = 0.004
$ws1.Range("C15").Value = 0.949
$ws1.Range("D15").Value = 0.905
$ws1.Range("B16").Value = 0.007
$ws1.Range("C16").Value = 0.987
$ws1.Range("D16").Value = 0.909
$ws1.Range("B17").Value = 0.005
$ws1.Range("C17").Value = 0.959
$ws1.Range("D17").Value = 0.905
$ws1.Range("B18").Value = 0.002
$ws1.Range("C18").Value = 0.924
$ws1.Range("D18").Value = 0.903
$ws1.Range("B19").Value = 0.006
$ws1.Range("C19").Value = 0.979
$ws1.Range("D19").Value = 0.913
$ws1.Range("B20").Value = -0.008
$ws1.Range("C20").Value = 0.803
$ws1.Range("D20").Value = 0.882
$ws1.Range("B21").Value = 0.007
$ws1.Range("C21").Value = 0.973
$ws1.Range("D21").Value = 0.897
$ws1.Range("B22").Value = 0.013
$ws1.Range("C22").Value = 1.081
$ws1.Range("D22").Value = 0.928
$ws1.Range("B23").Value = 0.025
$ws1.Range("C23").Value = 1.245
$ws1.Range("D23").Value = 0.928
$ws1.Range("B24").Value = 0.007
$ws1.Range("C24").Value = 0.994
$ws1.Range("D24").Value = 0.916
$ws1.Range("B25").Value = 0.01
$ws1.Range("C25").Value = 1.025
$ws1.Range("D25").Value = 0.912
$ws1.Range("B26").Value = -0.002
$ws1.Range("C26").Value = 0.893
$ws1.Range("D26").Value = 0.914
$ws1.Range("B27").Value = 0.001
$ws1.Range("C27").Value = 0.938
$ws1.Range("D27").Value = 0.927
$ws1.Range("B28").Value = 0.001
$ws1.Range("C28").Value = 0.94
$ws1.Range("B29").Value = 0.006
$ws1.Range("C29").Value = 0.978
$ws1.Range("D29").Value = 0.912
$ws1.Range("B30").Value = -0.01
$ws1.Range("C30").Value = 0.83
$ws1.Range("D30").Value = 0.935
$ws1.Range("B31").Value = 0.002
$ws1.Range("C31").Value = 0.967
$ws1.Range("D31").Value = 0.945
$ws1.Range("B32").Value = 0.02
$ws1.Range("C32").Value = 1.159
$ws1.Range("D32").Value = 0.917
$ws1.Range("B33").Value = 0.002
$ws1.Range("C33").Value = 0.948
$ws1.Range("D33").Value = 0.926
$ws1.Range("B34").Value = 0.026
$ws1.Range("C34").Value = 1.225
$ws1.Range("D34").Value = 0.901
$ws1.Range("B35").Value = 0.014
$ws1.Range("C35").Value = 1.079
$ws1.Range("D35").Value = 0.916
$ws1.Range("B36").Value = -0.009
$ws1.Range("C36").Value = 0.817
$ws1.Range("D36").Value = 0.908
$ws1.Range("B37").Value = 0.009
$ws1.Range("C37").Value = 1.009
$ws1.Range("D37").Value = 0.909
$ws1.Range("B38").Value = -0.006
$ws1.Range("C38").Value = 0.844
$ws1.Range("D38").Value = 0.905
$ws1.Range("B39").Value = 0.017
$ws1.Range("C39").Value = 1.131
$ws1.Range("D39").Value = 0.927
$ws1.Range("B40").Value = 0.002
$ws1.Range("C40").Value = 0.948
$ws1.Range("D40").Value = 0.926
$ws1.Range("B41").Value = 0.019
$ws1.Range("C41").Value = 1.142
$ws1.Range("D41").Value = 0.915
$ws1.Range("B42").Value = 0.01
$ws1.Range("C42").Value = 1.038
$ws1.Range("D42").Value = 0.923
$ws1.Range("B43").Value = -0.015
$ws1.Range("C43").Value = 0.77
$ws1.Range("D43").Value = 0.919
$ws1.Range("B44").Value = -0.007
$ws1.Range("C44").Value = 0.863
$ws1.Range("D44").Value = 0.937
$ws1.Range("B45").Value = -0.013
$ws1.Range("C45").Value = 0.789
$ws1.Range("D45").Value = 0.919
$ws1.Range("B46").Value = -0.004
$ws1.Range("C46").Value = 0.88
$ws1.Range("D46").Value = 0.922
$ws1.Range("B47").Value = 0.022
$ws1.Range("C47").Value = 1.168
$ws1.Range("D47").Value = 0.903
$ws1.Range("B48").Value = 0.005
$ws1.Range("C48").Value = 0.977
$ws1.Range("D48").Value = 0.922
$ws1.Range("B49").Value = -0.002
$ws1.Range("C49").Value = 0.904
$ws1.Range("D49").Value = 0.926
$ws1.Range("B50").Value = 0.016
$ws1.Range("C50").Value = 1.097
$ws1.Range("D50").Value = 0.909
$ws1.Range("B51").Value = 0.019
$ws1.Range("C51").Value = 1.131
$ws1.Range("D51").Value = 0.906
$ws1.Range("B52").Value = -0.01
$ws1.Range("C52").Value = 0.766
$ws1.Range("D52").Value = 0.862

$ws2.Range("B2").Value = -0.021
$ws2.Range("C2").Value = 0.72
$ws2.Range("D2").Value = 0.921
$ws2.Range("B3").Value = 0.014
$ws2.Range("C3").Value = 1.108
$ws2.Range("D3").Value = 0.94
$ws2.Range("B4").Value = -0.011
$ws2.Range("C4").Value = 0.82
$ws2.Range("D4").Value = 0.934
$ws2.Range("B5").Value = -0.007
$ws2.Range("C5").Value = 0.855
$ws2.Range("D5").Value = 0.928
$ws2.Range("B6").Value = 0.002
$ws2.Range("C6").Value = 0.958
$ws2.Range("D6").Value = 0.936
$ws2.Range("B7").Value = -0.012
$ws2.Range("C7").Value = 0.806
$ws2.Range("D7").Value = 0.928
$ws2.Range("B8").Value = -0.009
$ws2.Range("C8").Value = 0.823
$ws2.Range("D8").Value = 0.915
$ws2.Range("B9").Value = -0.002
$ws2.Range("C9").Value = 0.895
$ws2.Range("D9").Value = 0.916
$ws2.Range("B10").Value = -0.001
$ws2.Range("C10").Value = 0.909
$ws2.Range("D10").Value = 0.92
$ws2.Range("B11").Value = -0.016
$ws2.Range("C11").Value = 0.777
$ws2.Range("D11").Value = 0.94
$ws2.Range("B12").Value = -0.007
$ws2.Range("C12").Value = 0.847
$ws2.Range("D12").Value = 0.919
$ws2.Range("B13").Value = -0.006
$ws2.Range("C13").Value = 0.861
$ws2.Range("D13").Value = 0.924
$ws2.Range("B14").Value = -0.005
$ws2.Range("C14").Value = 0.869
$ws2.Range("D14").Value = 0.922
$ws2.Range("B15").Value = -0.016
$ws2.Range("C15").Value = 0.752
$ws2.Range("D15").Value = 0.907
$ws2.Range("B16").Value = -0.009
$ws2.Range("C16").Value = 0.821
$ws2.Range("D16").Value = 0.912
$ws2.Range("B17").Value = -0.009
$ws2.Range("C17").Value = 0.814
$ws2.Range("D17").Value = 0.904
$ws2.Range("B18").Value = -0.011
$ws2.Range("C18").Value = 0.8
$ws2.Range("D18").Value = 0.91
$ws2.Range("C19").Value = 0.834
$ws2.Range("D19").Value = 0.915
$ws2.Range("B20").Value = -0.02
$ws2.Range("C20").Value = 0.711
$ws2.Range("D20").Value = 0.899
$ws2.Range("B21").Value = -0.007
$ws2.Range("C21").Value = 0.833
$ws2.Range("D21").Value = 0.904
$ws2.Range("B22").Value = -0.006
$ws2.Range("C22").Value = 0.865
$ws2.Range("D22").Value = 0.928
$ws2.Range("B23").Value = -0.002
$ws2.Range("C23").Value = 0.909
$ws2.Range("D23").Value = 0.931
$ws2.Range("C24").Value = 0.827
$ws2.Range("D24").Value = 0.919
$ws2.Range("B25").Value = -0.003
$ws2.Range("C25").Value = 0.887
$ws2.Range("D25").Value = 0.919
$ws2.Range("B26").Value = -0.01
$ws2.Range("C26").Value = 0.815
$ws2.Range("D26").Value = 0.917
$ws2.Range("B27").Value = -0.008
$ws2.Range("C27").Value = 0.84
$ws2.Range("D27").Value = 0.923
$ws2.Range("B28").Value = -0.011
$ws2.Range("C28").Value = 0.81
$ws2.Range("D28").Value = 0.922
$ws2.Range("B29").Value = -0.009
$ws2.Range("C29").Value = 0.825
$ws2.Range("D29").Value = 0.916
$ws2.Range("B30").Value = -0.015
$ws2.Range("C30").Value = 0.786
$ws2.Range("D30").Value = 0.939
$ws2.Range("C31").Value = 0.904
$ws2.Range("D31").Value = 0.947
$ws2.Range("B32").Value = -0.002
$ws2.Range("C32").Value = 0.9
$ws2.Range("B33").Value = -0.011
$ws2.Range("C33").Value = 0.811
$ws2.Range("D33").Value = 0.922
$ws2.Range("B34").Value = -0.007
$ws2.Range("C34").Value = 0.83
$ws2.Range("D34").Value = 0.901
$ws2.Range("B35").Value = -0.0
$ws2.Range("C35").Value = 0.915
$ws2.Range("D35").Value = 0.915
$ws2.Range("B36").Value = -0.018
$ws2.Range("C36").Value = 0.741
$ws2.Range("D36").Value = 0.915
$ws2.Range("B37").Value = -0.008
$ws2.Range("C37").Value = 0.837
$ws2.Range("D37").Value = 0.919
$ws2.Range("B38").Value = -0.014
$ws2.Range("C38").Value = 0.764
$ws2.Range("D38").Value = 0.9
$ws2.Range("B39").Value = 0.001
$ws2.Range("C39").Value = 0.937
$ws2.Range("D39").Value = 0.926
$ws2.Range("B40").Value = -0.01
$ws2.Range("C40").Value = 0.821
$ws2.Range("D40").Value = 0.923
$ws2.Range("B41").Value = 0.002
$ws2.Range("C41").Value = 0.931
$ws2.Range("D41").Value = 0.909
$ws2.Range("B42").Value = -0.002
$ws2.Range("C42").Value = 0.9
$ws2.Range("D42").Value = 0.921
$ws2.Range("B43").Value = -0.023
$ws2.Range("C43").Value = 0.701
$ws2.Range("D43").Value = 0.919
$ws2.Range("B44").Value = -0.012
$ws2.Range("C44").Value = 0.809
$ws2.Range("D44").Value = 0.932
$ws2.Range("B45").Value = -0.019
$ws2.Range("C45").Value = 0.733
$ws2.Range("D45").Value = 0.917
$ws2.Range("B46").Value = -0.014
$ws2.Range("C46").Value = 0.778
$ws2.Range("D46").Value = 0.916
$ws2.Range("B47").Value = 0.001
$ws2.Range("C47").Value = 0.921
$ws2.Range("D47").Value = 0.91
$ws2.Range("B48").Value = -0.006
$ws2.Range("C48").Value = 0.859
$ws2.Range("D48").Value = 0.921
$ws2.Range("B49").Value = -0.014
$ws2.Range("C49").Value = 0.783
$ws2.Range("D49").Value = 0.923
$ws2.Range("B50").Value = -0.001
$ws2.Range("C50").Value = 0.891
$ws2.Range("D50").Value = 0.902
$ws2.Range("B51").Value = -0.009
$ws2.Range("C51").Value = 0.815
$ws2.Range("D51").Value = 0.905
$ws2.Range("B52").Value = -0.036
$ws2.Range("C52").Value = 0.563
$ws2.Range("D52").Value = 0.862

$ws3.Range("B2").Value = -0.027
$ws3.Range("C2").Value = 0.66
$ws3.Range("D2").Value = 907.56
$ws3.Range("E2").Value = 0.007
$ws3.Range("B3").Value = -0.041
$ws3.Range("C3").Value = 0.54
$ws3.Range("D3").Value = 7247.88
$ws3.Range("E3").Value = 0.058
$ws3.Range("B4").Value = -0.02
$ws3.Range("C4").Value = 0.726
$ws3.Range("D4").Value = 1386.17
$ws3.Range("E4").Value = 0.01
$ws3.Range("B5").Value = -0.016
$ws3.Range("C5").Value = 0.753
$ws3.Range("D5").Value = 1284.41
$ws3.Range("E5").Value = 0.011
$ws3.Range("B6").Value = -0.021
$ws3.Range("C6").Value = 0.702
$ws3.Range("D6").Value = 3593.61
$ws3.Range("E6").Value = 0.027
$ws3.Range("B7").Value = -0.025
$ws3.Range("C7").Value = 0.678
$ws3.Range("D7").Value = 1783.17
$ws3.Range("E7").Value = 0.014
$ws3.Range("B8").Value = -0.034
$ws3.Range("C8").Value = 0.588
$ws3.Range("D8").Value = 3873.09
$ws3.Range("E8").Value = 0.029
$ws3.Range("B9").Value = -0.029
$ws3.Range("C9").Value = 0.625
$ws3.Range("D9").Value = 3631.55
$ws3.Range("E9").Value = 0.031
$ws3.Range("B10").Value = -0.022
$ws3.Range("C10").Value = 0.684
$ws3.Range("D10").Value = 3107.85
$ws3.Range("E10").Value = 0.024
$ws3.Range("B11").Value = -0.019
$ws3.Range("C11").Value = 0.74
$ws3.Range("D11").Value = 645.42
$ws3.Range("E11").Value = 0.004
$ws3.Range("B12").Value = -0.025
$ws3.Range("C12").Value = 0.672
$ws3.Range("D12").Value = 2315.88
$ws3.Range("E12").Value = 0.02
$ws3.Range("B13").Value = -0.027
$ws3.Range("C13").Value = 0.649
$ws3.Range("D13").Value = 2927.96
$ws3.Range("E13").Value = 0.024
$ws3.Range("B14").Value = -0.03
$ws3.Range("C14").Value = 0.629
$ws3.Range("D14").Value = 3446.21
$ws3.Range("E14").Value = 0.028
$ws3.Range("B15").Value = -0.041
$ws3.Range("C15").Value = 0.539
$ws3.Range("D15").Value = 3515.02
$ws3.Range("E15").Value = 0.028
$ws3.Range("B16").Value = -0.029
$ws3.Range("C16").Value = 0.626
$ws3.Range("D16").Value = 2719.11
$ws3.Range("E16").Value = 0.023
$ws3.Range("B17").Value = -0.023
$ws3.Range("C17").Value = 0.685
$ws3.Range("D17").Value = 1540.31
$ws3.Range("E17").Value = 0.013
$ws3.Range("B18").Value = -0.027
$ws3.Range("C18").Value = 0.643
$ws3.Range("D18").Value = 2288.54
$ws3.Range("E18").Value = 0.018
$ws3.Range("B19").Value = -0.023
$ws3.Range("C19").Value = 0.686
$ws3.Range("D19").Value = 1819.81
$ws3.Range("E19").Value = 0.016
$ws3.Range("B20").Value = -0.035
$ws3.Range("C20").Value = 0.578
$ws3.Range("D20").Value = 2212.68
$ws3.Range("E20").Value = 0.018
$ws3.Range("B21").Value = -0.037
$ws3.Range("C21").Value = 0.555
$ws3.Range("D21").Value = 3958.54
$ws3.Range("E21").Value = 0.034
$ws3.Range("B22").Value = -0.028
$ws3.Range("C22").Value = 0.645
$ws3.Range("D22").Value = 3006.67
$ws3.Range("E22").Value = 0.025
$ws3.Range("B23").Value = -0.036
$ws3.Range("C23").Value = 0.577
$ws3.Range("D23").Value = 4849.68
$ws3.Range("E23").Value = 0.038
$ws3.Range("B24").Value = -0.03
$ws3.Range("C24").Value = 0.627
$ws3.Range("D24").Value = 2862.95
$ws3.Range("E24").Value = 0.023
$ws3.Range("B25").Value = -0.03
$ws3.Range("C25").Value = 0.617
$ws3.Range("D25").Value = 3553.18
$ws3.Range("E25").Value = 0.03
$ws3.Range("B26").Value = -0.024
$ws3.Range("C26").Value = 0.675
$ws3.Range("D26").Value = 2049.18
$ws3.Range("E26").Value = 0.016
$ws3.Range("B27").Value = -0.017
$ws3.Range("C27").Value = 0.748
$ws3.Range("D27").Value = 1186.41
$ws3.Range("E27").Value = 0.009
$ws3.Range("B28").Value = -0.036
$ws3.Range("C28").Value = 0.584
$ws3.Range("D28").Value = 3191.38
$ws3.Range("E28").Value = 0.025
$ws3.Range("B29").Value = -0.031
$ws3.Range("C29").Value = 0.617
$ws3.Range("D29").Value = 2795.17
$ws3.Range("E29").Value = 0.024
$ws3.Range("B30").Value = -0.018
$ws3.Range("C30").Value = 0.759
$ws3.Range("D30").Value = 456.86
$ws3.Range("E30").Value = 0.003
$ws3.Range("B31").Value = -0.015
$ws3.Range("C31").Value = 0.775
$ws3.Range("D31").Value = 1651.0
$ws3.Range("E31").Value = 0.014
$ws3.Range("B32").Value = -0.036
$ws3.Range("C32").Value = 0.572
$ws3.Range("D32").Value = 4650.02
$ws3.Range("E32").Value = 0.038
$ws3.Range("B33").Value = -0.032
$ws3.Range("C33").Value = 0.618
$ws3.Range("D33").Value = 2658.52
$ws3.Range("E33").Value = 0.022
$ws3.Range("B34").Value = -0.042
$ws3.Range("C34").Value = 0.514
$ws3.Range("D34").Value = 5626.01
$ws3.Range("E34").Value = 0.042
$ws3.Range("B35").Value = -0.017
$ws3.Range("C35").Value = 0.735
$ws3.Range("D35").Value = 2146.52
$ws3.Range("E35").Value = 0.018
$ws3.Range("B36").Value = -0.029
$ws3.Range("C36").Value = 0.639
$ws3.Range("D36").Value = 1402.51
$ws3.Range("E36").Value = 0.012
$ws3.Range("B37").Value = -0.03
$ws3.Range("C37").Value = 0.624
$ws3.Range("D37").Value = 2825.32
$ws3.Range("E37").Value = 0.025
$ws3.Range("B38").Value = -0.02
$ws3.Range("C38").Value = 0.704
$ws3.Range("D38").Value = 825.83
$ws3.Range("E38").Value = 0.007
$ws3.Range("B39").Value = -0.017
$ws3.Range("C39").Value = 0.742
$ws3.Range("D39").Value = 2268.2
$ws3.Range("E39").Value = 0.019
$ws3.Range("B40").Value = -0.026
$ws3.Range("C40").Value = 0.667
$ws3.Range("D40").Value = 2171.13
$ws3.Range("E40").Value = 0.018
$ws3.Range("B41").Value = -0.028
$ws3.Range("C41").Value = 0.624
$ws3.Range("D41").Value = 3977.99
$ws3.Range("E41").Value = 0.033
$ws3.Range("B42").Value = -0.019
$ws3.Range("C42").Value = 0.713
$ws3.Range("D42").Value = 2411.23
$ws3.Range("E42").Value = 0.02
$ws3.Range("B43").Value = -0.03
$ws3.Range("C43").Value = 0.641
$ws3.Range("D43").Value = 939.64
$ws3.Range("E43").Value = 0.007
$ws3.Range("B44").Value = -0.017
$ws3.Range("C44").Value = 0.763
$ws3.Range("D44").Value = 636.54
$ws3.Range("E44").Value = 0.005
$ws3.Range("B45").Value = -0.029
$ws3.Range("C45").Value = 0.641
$ws3.Range("D45").Value = 1479.18
$ws3.Range("E45").Value = 0.011
$ws3.Range("B46").Value = -0.028
$ws3.Range("C46").Value = 0.648
$ws3.Range("D46").Value = 1869.81
$ws3.Range("E46").Value = 0.016
$ws3.Range("B47").Value = -0.033
$ws3.Range("C47").Value = 0.583
$ws3.Range("D47").Value = 4606.46
$ws3.Range("E47").Value = 0.038
$ws3.Range("B48").Value = -0.023
$ws3.Range("C48").Value = 0.685
$ws3.Range("D48").Value = 2264.74
$ws3.Range("E48").Value = 0.019
$ws3.Range("B49").Value = -0.03
$ws3.Range("C49").Value = 0.639
$ws3.Range("D49").Value = 2182.0
$ws3.Range("E49").Value = 0.017
$ws3.Range("B50").Value = -0.039
$ws3.Range("C50").Value = 0.54
$ws3.Range("D50").Value = 5143.42
$ws3.Range("E50").Value = 0.042
$ws3.Range("B51").Value = -0.035
$ws3.Range("C51").Value = 0.576
$ws3.Range("D51").Value = 3347.89
$ws3.Range("E51").Value = 0.028
$ws3.Range("B52").Value = -0.057
$ws3.Range("C52").Value = 0.422
$ws3.Range("D52").Value = 3540.6
$ws3.Range("E52").Value = 0.026

Write-Host "Applied all changes"
